# Update countries & provincias Spain
# Applies the "paises.xlsx" data refresh:
#   - Adds "Liberia" as a newly-tracked country (placed right after "Martinica")
#     with its own case counts, which pushes Guadalupe/Birmania/Gibraltar down
#     one row each (their data is unchanged) while Brunei (and everything after)
#     keeps its original row.
#   - Refreshes case-count figures for several existing countries.
#   - Bumps the "Datos actualizados..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert "Liberia" as a new row right after "Martinica" (row 136), without
#    changing the overall row count: first drop the old Liberia row (140),
#    then insert a fresh blank row at 137 and fill it in. Net effect: every
#    row from Guadalupe (old row 137) through Gibraltar (old row 139) simply
#    slides down to 138-140 untouched, Brunei stays put at row 141, and the
#    sheet's last row remains 218.
# ---------------------------------------------------------------------------
$ws.Rows(140).Delete()
$ws.Rows(137).Insert()

$ws.Range("A137").Value = "Liberia"
$ws.Range("B137").Value = 152
$ws.Range("C137").Value = 11
$ws.Range("D137").Value = 45
$ws.Range("E137").Value = 89
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 2
$ws.Range("H137").Value = 18

# ---------------------------------------------------------------------------
# 2) Refresh figures for countries whose counts changed in this update.
# ---------------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1099275
$ws.Range("C4").Value = 4252
$ws.Range("D4").Value = 156089
$ws.Range("E4").Value = 879214
$ws.Range("F4").Value = 15247
$ws.Range("G4").Value = 116
$ws.Range("H4").Value = 63972

# Row 5 - España
$ws.Range("E5").Value = 76831
$ws.Range("F5").Value = 2500
$ws.Range("G5").Value = 281
$ws.Range("H5").Value = 24824

# Row 9 - Alemania
$ws.Range("B9").Value = 163331
$ws.Range("C9").Value = 322
$ws.Range("E9").Value = 29799

# Row 20 - Suiza
$ws.Range("E20").Value = 4556
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 1749

# Row 122 - Kenia
$ws.Range("B122").Value = 411
$ws.Range("C122").Value = 15
$ws.Range("D122").Value = 150
$ws.Range("E122").Value = 240
$ws.Range("G122").Value = 4
$ws.Range("H122").Value = 21

# Row 124 - Venezuela
$ws.Range("E124").Value = 181
$ws.Range("H124").Value = 10

# ---------------------------------------------------------------------------
# 3) Bump the "last updated" timestamp banner.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 16:22"
